$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 29.20950566666667
$ws.Range("H2").Value = 87.628517
$ws.Range("I2").Value = 0.01829497698069002
$ws.Range("J2").Value = 0.01840828041918582
$ws.Range("M2").Value = 1.009860666666667
$ws.Range("N2").Value = 3.029582
$ws.Range("O2").Value = 0.01353413605720072
$ws.Range("P2").Value = 0.01542521070970148
$ws.Range("Q2").Value = 29.49753086554378
$ws.Range("R2").Value = 265.477777789894
$ws.Range("S2").Value = 0.000247606707620014
$ws.Range("T2").Value = 0.0002839516042692132

$ws.Range("G3").Value = 29.20950566666667
$ws.Range("H3").Value = 87.628517
$ws.Range("I3").Value = 0.01829497698069002
$ws.Range("J3").Value = 0.01840828041918582
$ws.Range("O3").Value = 0.6185519418990597
$ws.Range("P3").Value = 0.704979911415303
$ws.Range("Q3").Value = 1348.128533730977
$ws.Range("R3").Value = 12133.15680357879
$ws.Range("S3").Value = 0.01131639353840441
$ws.Range("T3").Value = 0.01297746789922568

$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("G4").Value = 29.20950566666667
$ws.Range("H4").Value = 87.628517
$ws.Range("I4").Value = 0.01829497698069002
$ws.Range("J4").Value = 0.01840828041918582
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.009315666666666667
$ws.Range("N4").Value = 0.027947
$ws.Range("O4").Value = 0.0001248484115599408
$ws.Range("P4").Value = 0.000142293017222847
$ws.Range("Q4").Value = 0.2721060182887778
$ws.Range("R4").Value = 2.448954164599
$ws.Range("S4").Value = 0.000002284098815564831
$ws.Range("T4").Value = 0.000002619369762730205

$ws.Range("D5").Value = "MuSCs"
$ws.Range("G5").Value = 29.20950566666667
$ws.Range("H5").Value = 87.628517
$ws.Range("I5").Value = 0.01829497698069002
$ws.Range("J5").Value = 0.01840828041918582
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 27.4428835
$ws.Range("N5").Value = 54.885767
$ws.Range("O5").Value = 0.3677890736321797
$ws.Range("P5").Value = 0.2794525848577725
$ws.Range("Q5").Value = 801.5930611029232
$ws.Range("R5").Value = 4809.558366617539
$ws.Range("S5").Value = 0.006728692635850034
$ws.Range("T5").Value = 0.005144241545928198

$ws.Range("I6").Value = 0.913374480506715
$ws.Range("J6").Value = 0.9190311407684336
$ws.Range("M6").Value = 1.009860666666667
$ws.Range("N6").Value = 3.029582
$ws.Range("O6").Value = 0.01353413605720072
$ws.Range("P6").Value = 0.01542521070970148
$ws.Range("Q6").Value = 1472.660608372664
$ws.Range("R6").Value = 13253.94547535397
$ws.Range("S6").Value = 0.01236173449035291
$ws.Range("T6").Value = 0.01417624899513041

$ws.Range("I7").Value = 0.913374480506715
$ws.Range("J7").Value = 0.9190311407684336
$ws.Range("O7").Value = 0.6185519418990597
$ws.Range("P7").Value = 0.704979911415303
$ws.Range("S7").Value = 0.5649695585984734
$ws.Range("T7").Value = 0.6478984922068353

$ws.Range("D8").Value = "Inflammatory-Mac"
$ws.Range("I8").Value = 0.913374480506715
$ws.Range("J8").Value = 0.9190311407684336
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.009315666666666667
$ws.Range("N8").Value = 0.027947
$ws.Range("O8").Value = 0.0001248484115599408
$ws.Range("P8").Value = 0.000142293017222847
$ws.Range("Q8").Value = 13.58485956880878
$ws.Range("R8").Value = 122.263736119279
$ws.Range("S8").Value = 0.0001140333530506495
$ws.Range("T8").Value = 0.0001307717139416954

$ws.Range("D9").Value = "MuSCs"
$ws.Range("I9").Value = 0.913374480506715
$ws.Range("J9").Value = 0.9190311407684336
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 27.4428835
$ws.Range("N9").Value = 54.885767
$ws.Range("O9").Value = 0.3677890736321797
$ws.Range("P9").Value = 0.2794525848577725
$ws.Range("Q9").Value = 40019.4351999155
$ws.Range("R9").Value = 240116.611199493
$ws.Range("S9").Value = 0.335929154064838
$ws.Range("T9").Value = 0.2568256278525262

$ws.Range("G10").Value = 57.98602933333333
$ws.Range("H10").Value = 173.958088
$ws.Range("I10").Value = 0.03631876156896331
$ws.Range("J10").Value = 0.03654368891224535
$ws.Range("M10").Value = 1.009860666666667
$ws.Range("N10").Value = 3.029582
$ws.Range("O10").Value = 0.01353413605720072
$ws.Range("P10").Value = 0.01542521070970148
$ws.Range("Q10").Value = 58.55781023991289
$ws.Range("R10").Value = 527.020292159216
$ws.Range("S10").Value = 0.0004915430605033822
$ws.Range("T10").Value = 0.0005636941015811661

$ws.Range("G11").Value = 57.98602933333333
$ws.Range("H11").Value = 173.958088
$ws.Range("I11").Value = 0.03631876156896331
$ws.Range("J11").Value = 0.03654368891224535
$ws.Range("O11").Value = 0.6185519418990597
$ws.Range("P11").Value = 0.704979911415303
$ws.Range("Q11").Value = 2676.273319860979
$ws.Range("R11").Value = 24086.45987874881
$ws.Range("S11").Value = 0.0224650404958512
$ws.Range("T11").Value = 0.02576256657214311

$ws.Range("D12").Value = "Inflammatory-Mac"
$ws.Range("G12").Value = 57.98602933333333
$ws.Range("H12").Value = 173.958088
$ws.Range("I12").Value = 0.03631876156896331
$ws.Range("J12").Value = 0.03654368891224535
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.009315666666666667
$ws.Range("N12").Value = 0.027947
$ws.Range("O12").Value = 0.0001248484115599408
$ws.Range("P12").Value = 0.000142293017222847
$ws.Range("Q12").Value = 0.5401785205928888
$ws.Range("R12").Value = 4.861606685336
$ws.Range("S12").Value = 0.000004534339691709293
$ws.Range("T12").Value = 0.00000519991175577649

$ws.Range("D13").Value = "MuSCs"
$ws.Range("G13").Value = 57.98602933333333
$ws.Range("H13").Value = 173.958088
$ws.Range("I13").Value = 0.03631876156896331
$ws.Range("J13").Value = 0.03654368891224535
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 27.4428835
$ws.Range("N13").Value = 54.885767
$ws.Range("O13").Value = 0.3677890736321797
$ws.Range("P13").Value = 0.2794525848577725
$ws.Range("Q13").Value = 1591.303847622249
$ws.Range("R13").Value = 9547.823085733495
$ws.Range("S13").Value = 0.01335764367291702
$ws.Range("T13").Value = 0.01021222832676528

$ws.Range("G14").Value = 29.481085
$ws.Range("H14").Value = 58.96217
$ws.Range("I14").Value = 0.01846507700595112
$ws.Range("J14").Value = 0.01238628926567028
$ws.Range("M14").Value = 1.009860666666667
$ws.Range("N14").Value = 3.029582
$ws.Range("O14").Value = 0.01353413605720072
$ws.Range("P14").Value = 0.01542521070970148
$ws.Range("Q14").Value = 29.77178815215667
$ws.Range("R14").Value = 178.63072891294
$ws.Range("S14").Value = 0.0002499088645052311
$ws.Range("T14").Value = 0.0001910611218342776

$ws.Range("G15").Value = 29.481085
$ws.Range("H15").Value = 58.96217
$ws.Range("I15").Value = 0.01846507700595112
$ws.Range("J15").Value = 0.01238628926567028
$ws.Range("O15").Value = 0.6185519418990597
$ws.Range("P15").Value = 0.704979911415303
$ws.Range("Q15").Value = 1360.662941283657
$ws.Range("R15").Value = 8163.977647701939
$ws.Range("S15").Value = 0.01142160923934674
$ws.Range("T15").Value = 0.00873208510927655

$ws.Range("D16").Value = "Inflammatory-Mac"
$ws.Range("G16").Value = 29.481085
$ws.Range("H16").Value = 58.96217
$ws.Range("I16").Value = 0.01846507700595112
$ws.Range("J16").Value = 0.01238628926567028
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.3333333333333333
$ws.Range("M16").Value = 0.009315666666666667
$ws.Range("N16").Value = 0.027947
$ws.Range("O16").Value = 0.0001248484115599408
$ws.Range("P16").Value = 0.000142293017222847
$ws.Range("Q16").Value = 0.2746359608316667
$ws.Range("R16").Value = 1.64781576499
$ws.Range("S16").Value = 0.000002305335533524985
$ws.Range("T16").Value = 0.000001762482471807186

$ws.Range("D17").Value = "MuSCs"
$ws.Range("G17").Value = 29.481085
$ws.Range("H17").Value = 58.96217
$ws.Range("I17").Value = 0.01846507700595112
$ws.Range("J17").Value = 0.01238628926567028
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 27.4428835
$ws.Range("N17").Value = 54.885767
$ws.Range("O17").Value = 0.3677890736321797
$ws.Range("P17").Value = 0.2794525848577725
$ws.Range("Q17").Value = 809.0459811085975
$ws.Range("R17").Value = 3236.18392443439
$ws.Range("S17").Value = 0.006791253566565626
$ws.Range("T17").Value = 0.00346138055208764

$ws.Range("G18").Value = 21.628479
$ws.Range("H18").Value = 64.885437
$ws.Range("I18").Value = 0.01354670393768061
$ws.Range("J18").Value = 0.01363060063446486
$ws.Range("M18").Value = 1.009860666666667
$ws.Range("N18").Value = 3.029582
$ws.Range("O18").Value = 0.01353413605720072
$ws.Range("P18").Value = 0.01542521070970148
$ws.Range("Q18").Value = 21.841750221926
$ws.Range("R18").Value = 196.575751997334
$ws.Range("S18").Value = 0.0001833429342191862
$ws.Range("T18").Value = 0.0002102548868864112

$ws.Range("G19").Value = 21.628479
$ws.Range("H19").Value = 64.885437
$ws.Range("I19").Value = 0.01354670393768061
$ws.Range("J19").Value = 0.01363060063446486
$ws.Range("O19").Value = 0.6185519418990597
$ws.Range("P19").Value = 0.704979911415303
$ws.Range("Q19").Value = 998.2356433500258
$ws.Range("R19").Value = 8984.120790150233
$ws.Range("S19").Value = 0.00837934002698398
$ws.Range("T19").Value = 0.009609299627822412

$ws.Range("D20").Value = "Inflammatory-Mac"
$ws.Range("G20").Value = 21.628479
$ws.Range("H20").Value = 64.885437
$ws.Range("I20").Value = 0.01354670393768061
$ws.Range("J20").Value = 0.01363060063446486
$ws.Range("K20").Value = 1
$ws.Range("L20").Value = 0.3333333333333333
$ws.Range("M20").Value = 0.009315666666666667
$ws.Range("N20").Value = 0.027947
$ws.Range("O20").Value = 0.0001248484115599408
$ws.Range("P20").Value = 0.000142293017222847
$ws.Range("Q20").Value = 0.201483700871
$ws.Range("R20").Value = 1.813353307839
$ws.Range("S20").Value = 0.00000169128446849222
$ws.Range("T20").Value = 0.000001939539290837658

$ws.Range("D21").Value = "MuSCs"
$ws.Range("G21").Value = 21.628479
$ws.Range("H21").Value = 64.885437
$ws.Range("I21").Value = 0.01354670393768061
$ws.Range("J21").Value = 0.01363060063446486
$ws.Range("K21").Value = 2
$ws.Range("L21").Value = 1
$ws.Range("M21").Value = 27.4428835
$ws.Range("N21").Value = 54.885767
$ws.Range("O21").Value = 0.3677890736321797
$ws.Range("P21").Value = 0.2794525848577725
$ws.Range("Q21").Value = 593.5478294791965
$ws.Range("R21").Value = 3561.286976875179
$ws.Range("S21").Value = 0.004982329692008952
$ws.Range("T21").Value = 0.0038091065804652
